$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (row 3) of data for IP 122.180.21.165, filled column-by-column
# (D, E, F, I) to mirror the order new shared strings were introduced in
# the source edit, then the remaining plain cells.

# Column D - Link (text + hyperlink, same visual style as D2)
$ws.Range("D3").Value = "https://www.virustotal.com/gui/ip-address/122.180.21.165/detection"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.virustotal.com/gui/ip-address/122.180.21.165/detection") | Out-Null
$ws.Range("D3").Style = $ws.Range("D2").Style

# Column E - last_analysis_stats
$ws.Range("E3").Value = "{'harmless': 57, 'malicious': 11, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"

# Column F - Country
$ws.Range("F3").Value = "India"

# Column I - AS_Owner
$ws.Range("I3").Value = "Bharti Airtel Ltd., Telemedia Services"

# Remaining cells for the new row
$ws.Range("C3").Value = "Malicious"
$ws.Range("G3").Value = 45328.50989583333
$ws.Range("H3").Value = 45340.69359953704
$ws.Range("G3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wb.Save()
